# "Generate Report for Handoff"
# The b.md file has just been handed off for localization (zh-cn and de-de),
# so its status rows move from "Handed back: in sync with en-US" to
# "Ready for handoff", with fresh handoff file names / timestamps, and the
# de-de row now also carries a "stale handback" error detail message.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$overviewDate = "2016-08-29 02:37:06"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/12429caebe9e81c964138f81e951b33e26c84e73/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bd6cf20f51b0992b423309cd14926efbf52e29ae/e2e/b.md."

# --- Overview sheet: row 3 is b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $overviewDate

# --- zh-cn sheet: row 3 is b.md ---
# NOTE: "Content Duplicate" (F) is a text column holding the literal words
# "True"/"False" (shared-string "s" cells in the original file), not real
# booleans, so a leading apostrophe is used to force text entry instead of
# letting the host auto-coerce "False" into a Boolean cell.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $status
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("F3").ClearFormats()
$wsZh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-29 02:36:58"
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is b.md ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $status
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("F3").ClearFormats()
$wsDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("H3").Value = $overviewDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.17
